# Handback status report: replace the in-flight file entry with its final
# name/hashes and append a second, newly-handed-back file as a new row on
# every sheet (Overview, zh-cn, de-de). Mirrors the "Generate Report for
# Handback" commit.

$wb = $excel.ActiveWorkbook

$oldId    = "ccb36728-a1a9-47cf-85bf-806625b840a0"
$newId1   = "4032afde-8e5a-4dd2-a7f0-f64482d5e798"
$newId2   = "961d70ea-fabc-4d5e-9a78-48ffc12531d6"

$hashZh1  = "75ee24a80cc4ad9e7f13a0a97fb1ea756794ef95"
$hashZh2  = "6654b063829a5948dac2225325525dd1ca8915b1"

$repoMain = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70990f5c9c6b5bfe6ca2d7a9818aa8f15bcf900b/e2e"
$repoZh   = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/953e18a8ed7785aba1fc043b6c1e345770382bb1/e2e"
$repoDe   = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4738f0411da4e396710bfc96a7942f0ff0360807/e2e"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Row 2 already holds the first file's entry - rename it to its final
# generated id and bump the generate-date to the real handback time.
$ws.Cells.Item(2, 1).Value = "$newId1.md"
$ws.Hyperlinks.Add($ws.Range("B2"), "$repoMain/$newId1.md", [Type]::Missing, [Type]::Missing, "e2e\$newId1.md") | Out-Null
$ws.Cells.Item(2, 7).Value = "2016-08-22 19:05:57"

# Row 3: newly handed-back second file.
$ws.Cells.Item(3, 1).Value = "$newId2.md"
$ws.Hyperlinks.Add($ws.Range("B3"), "$repoMain/$newId2.md", [Type]::Missing, [Type]::Missing, "e2e\$newId2.md") | Out-Null
$ws.Cells.Item(3, 3).Value = ".md"
$ws.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(3, 7).Value = "2016-08-22 19:05:57"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Row 2: update to the final generated id / hash / handoff+handback times.
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoMain/$newId1.md", [Type]::Missing, [Type]::Missing, "$newId1.md") | Out-Null
$ws.Cells.Item(2, 7).Value = "$newId1.$hashZh1.zh-cn.xlf"
$ws.Cells.Item(2, 8).Value = "2016-08-22 19:05:51"
$ws.Hyperlinks.Add($ws.Range("I2"), "$repoMain/$newId1.md", [Type]::Missing, [Type]::Missing, "$newId1.md") | Out-Null
$ws.Cells.Item(2, 10).Value = "$newId1.$hashZh1.zh-cn.xlf"
$ws.Cells.Item(2, 11).Value = "2016-08-22 19:06:24"

# Row 3: second file's zh-cn row.
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoMain/$newId2.md", [Type]::Missing, [Type]::Missing, "$newId2.md") | Out-Null
$ws.Cells.Item(3, 2).Value = ".md"
$ws.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(3, 4).Value = "e2e"
$ws.Cells.Item(3, 5).Value = "ht"
$ws.Cells.Item(3, 6).Value = "True"
$ws.Cells.Item(3, 7).Value = "$newId2.$hashZh2.zh-cn.xlf"
$ws.Cells.Item(3, 8).Value = "2016-08-22 19:06:32"
$ws.Hyperlinks.Add($ws.Range("I3"), "$repoMain/$newId2.md", [Type]::Missing, [Type]::Missing, "$newId2.md") | Out-Null
$ws.Cells.Item(3, 10).Value = "$newId2.$hashZh2.zh-cn.xlf"
$ws.Cells.Item(3, 11).Value = "2016-08-22 19:06:24"
$ws.Cells.Item(3, 13).Value = "True"
$ws.Cells.Item(3, 15).Value = "False"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

# Row 2: update to the final generated id / hash / handoff+handback times.
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoMain/$newId1.md", [Type]::Missing, [Type]::Missing, "$newId1.md") | Out-Null
$ws.Cells.Item(2, 7).Value = "$newId1.$hashZh1.de-de.xlf"
$ws.Cells.Item(2, 8).Value = "2016-08-22 19:05:57"
$ws.Hyperlinks.Add($ws.Range("I2"), "$repoMain/$newId1.md", [Type]::Missing, [Type]::Missing, "$newId1.md") | Out-Null
$ws.Cells.Item(2, 10).Value = "$newId1.$hashZh1.de-de.xlf"
$ws.Cells.Item(2, 11).Value = "2016-08-22 19:06:32"

# Row 3: second file's de-de row.
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoMain/$newId2.md", [Type]::Missing, [Type]::Missing, "$newId2.md") | Out-Null
$ws.Cells.Item(3, 2).Value = ".md"
$ws.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(3, 4).Value = "e2e"
$ws.Cells.Item(3, 5).Value = "ht"
$ws.Cells.Item(3, 6).Value = "True"
$ws.Cells.Item(3, 7).Value = "$newId2.$hashZh2.de-de.xlf"
$ws.Cells.Item(3, 8).Value = "2016-08-22 19:05:57"
$ws.Hyperlinks.Add($ws.Range("I3"), "$repoMain/$newId2.md", [Type]::Missing, [Type]::Missing, "$newId2.md") | Out-Null
$ws.Cells.Item(3, 10).Value = "$newId2.$hashZh2.de-de.xlf"
$ws.Cells.Item(3, 11).Value = "2016-08-22 19:06:32"
$ws.Cells.Item(3, 13).Value = "True"
$ws.Cells.Item(3, 15).Value = "False"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

Write-Host "Handback report regenerated."
